$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "draft data" task rows (NANNY / UNLIMITED CHILD CARE / SYP / JUNIOR /
# CLEANUP sections, previously rows 58-77) need to shift down three rows so
# two new notes plus an extra spacer row fit above them (new rows 58-60 are
# blank placeholders, matching the original's two blank rows pattern).
[void]$ws.Rows("58:60").Insert()

# Add the two new draft notes into the freed-up blank rows.
$ws.Range("D56").Value = "Need to add Ifee to transactions"
$ws.Range("D57").Value = "Add dropdown for Sales Reps to main page"

# Restore the saved selection/view state (now pointing at the newly added row).
[void]$ws.Range("D57").Select()
